# Adds a new FAQ question/answer block at the end of the document:
#   "¿Qué debo escribir cuándo el enunciado me dice "se dispone"?"
# with two bullet-level answer paragraphs, matching the existing
# question/answer list pattern used throughout the document.

$d = $word.ActiveDocument

# Start from the end of the last existing paragraph and create a fresh
# empty paragraph to receive the new content (keeps the existing last
# paragraph's text/formatting untouched).
$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

# --- Paragraph 1: the bolded question, list level 0 ("ilvl 0") ---
$p1 = $d.Paragraphs.Last.Range
$p1Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
      <w:lang w:val="es-AR"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>&#191;Qu&#233; debo escribir cu&#225;ndo el enunciado me dice &#8220;se dispone&#8221;?</w:t>
  </w:r>
</w:p>
'@
$p1.InsertXML($p1Xml)

# --- Paragraph 2: first answer paragraph, list level 1 ("ilvl 1") ---
$lastPara2 = $d.Paragraphs.Last
$endRange2 = $lastPara2.Range
$endRange2.Collapse(0)
$endRange2.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last.Range
$p2Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
      <w:lang w:val="es-AR"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t xml:space="preserve">Lo &#250;nico que hay que hacer cuando ya se dispone de una estructura es escribir el </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>type</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t xml:space="preserve"> completo y en el programa hacer un llamado </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>ej</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>: &#8220;</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>cargarLista</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>(</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>L);&#8221; y nada m&#225;s. No hay que implementar ning&#250;n modulo ni nada.</w:t>
  </w:r>
</w:p>
'@
$p2.InsertXML($p2Xml)

# --- Paragraph 3: second answer paragraph, list level 1 ("ilvl 1") ---
$lastPara3 = $d.Paragraphs.Last
$endRange3 = $lastPara3.Range
$endRange3.Collapse(0)
$endRange3.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last.Range
$p3Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
      <w:lang w:val="es-AR"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-AR"/>
    </w:rPr>
    <w:t>Recordar poner un comentario al lado del llamado de la carga diciendo &#8220;se dispone&#8221;.</w:t>
  </w:r>
</w:p>
'@
$p3.InsertXML($p3Xml)
